$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# The weather-log "runner" now scrapes a different set of cities (selenium
# sheet changed runner) -- replace the old Places/Coordinates/timeStamp table
# with the new Places/Temperature/Temperature1/Temperature2 table.
# ---------------------------------------------------------------------------

# Wipe the old table body (keeps header row formatting in place).
$ws.Range("A1:C8").ClearContents() | Out-Null

# Headers
$ws.Cells.Item(1, 1).Value = "Places"
$ws.Cells.Item(1, 2).Value = "Temperature"

# Column A - place names (rows 2-11 first, matching the scrape order)
$colA1 = @("jammu", "moga", "hisar", "palwal", "agra", "udaipur", "ajmer", "surat", "chiplun", "karwar")
for ($i = 0; $i -lt $colA1.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA1[$i]
}
$ws.Cells.Item(13, 1).Value = "kannur"
$ws.Cells.Item(14, 1).Value = "kozhikode"

# Column B - first temperature reading (only available for the first 10 cities)
$colB = @("33.9.", "33.02.", "35.99.", "31.2.", "34.15.", "25.5.", "23.42.", "29.99.", "28.22.", "27.11.")
for ($i = 0; $i -lt $colB.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

# New headers for the extra temperature readings
$ws.Cells.Item(1, 3).Value = "Temperature1"
$ws.Cells.Item(1, 4).Value = "Temperature2"

# Row 12 place name added afterwards
$ws.Cells.Item(12, 1).Value = "kasaragod"

# Column C - second temperature reading (all 13 cities)
$colC = @("28.17.", "29.84.", "29.88.", "30.32.", "28.89.", "22.56.", "23.9.", "27.02.", "24.99.", "25.75.", "26.59.", "25.37.", "24.65.")
for ($i = 0; $i -lt $colC.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

# Match the yellow header-row highlight on the new header cells.
$ws.Cells.Item(1, 3).Interior.Color = 65535
$ws.Cells.Item(1, 4).Interior.Color = 65535

# Resize columns for the new content (column A keeps its existing width).
$ws.Columns.Item(2).ColumnWidth = 12.0
$ws.Columns.Item(3).ColumnWidth = 13.0
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668

# Selection ends on the newly-added row (A12).
$ws.Range("A12").Select() | Out-Null
